# Applies the "Add files via upload" edit to inverter_charac.xlsx:
#  - Rename a few header labels in row 1 (units/casing clarifications)
#  - Fill in previously-blank W(Width) / Noise-Margin measurement columns
#    for data rows 2-16
#  - Update the view (selection / zoom) to match the author's last save

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) text updates ----------------------------------
$ws.Range("B1").Value = "W (Width)pmos(nm)"
$ws.Range("E1").Value = "W (Width)nmos(nm)"
$ws.Range("N1").Value = "Noise Margin (NMH)"
$ws.Range("O1").Value = "Noise Margin (NML)"
$ws.Range("R1").Value = "NFIN"

# ---- Data rows 2-16: fill in W(Width)pmos/nmos and Noise Margin values -

# N2 ends up with no explicit cell style (reverts to the workbook default),
# unlike its neighbours which keep style index 1 - clear its format first.
$ws.Range("N2").Style = "Normal"

$ws.Range("B2").Value = 987
$ws.Range("D2").Value = 141
$ws.Range("E2").Value = 987
$ws.Range("G2").Value = 141
$ws.Range("N2").Value = 0.19400000000000001
$ws.Range("O2").Value = 0.17599999999999999

$ws.Range("B3").Value = 1269
$ws.Range("D3").Value = 181
$ws.Range("E3").Value = 987
# G3 did not exist yet and is created using the "thin side borders only"
# format (the same one already used by N14) rather than the normal style.
$ws.Range("N14").Copy()
$ws.Range("G3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G3").Value = 141
$ws.Range("N3").Value = 0.16500000000000001
$ws.Range("O3").Value = 0.19800000000000001

$ws.Range("B4").Value = 705
$ws.Range("D4").Value = 100.71
$ws.Range("E4").Value = 987
$ws.Range("G4").Value = 141
$ws.Range("N4").Value = 0.23499999999999999
$ws.Range("O4").Value = 0.13300000000000001

$ws.Range("B5").Value = 1128
$ws.Range("D5").Value = 161.114
$ws.Range("E5").Value = 987
$ws.Range("G5").Value = 141
$ws.Range("N5").Value = 0.23599999999999999
$ws.Range("O5").Value = 0.14399999999999999

$ws.Range("B6").Value = 1057.5
$ws.Range("D6").Value = 151.07
$ws.Range("E6").Value = 987
$ws.Range("G6").Value = 141
$ws.Range("N6").Value = 0.20200000000000001
$ws.Range("O6").Value = 0.16300000000000001

$ws.Range("B7").Value = 1057.5
$ws.Range("D7").Value = 151.07
$ws.Range("E7").Value = 916.5
$ws.Range("G7").Value = 130.91999999999999
$ws.Range("N7").Value = 0.17899999999999999
$ws.Range("O7").Value = 0.186

$ws.Range("B8").Value = 1057.5
$ws.Range("D8").Value = 151.07
$ws.Range("E8").Value = 846
$ws.Range("G8").Value = 120.8
$ws.Range("N8").Value = 0.17
# O8 picks up the "thin side borders only" format instead of the normal one.
$ws.Range("N14").Copy()
$ws.Range("O8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O8").Value = 0.19400000000000001

$ws.Range("B9").Value = 987
$ws.Range("D9").Value = 141
$ws.Range("E9").Value = 916.5
$ws.Range("G9").Value = 130.91999999999999
$ws.Range("N9").Value = 0.187
# O9 also picks up the "thin side borders only" format.
$ws.Range("N14").Copy()
$ws.Range("O9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("O9").Value = 0.17699999999999999

$ws.Range("B10").Value = 916.5
$ws.Range("D10").Value = 130.9
# E10 also picks up the "thin side borders only" format.
$ws.Range("N14").Copy()
$ws.Range("E10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E10").Value = 846
$ws.Range("G10").Value = 120.8
$ws.Range("N10").Value = 0.185
$ws.Range("O10").Value = 0.17899999999999999

$ws.Range("B11").Value = 775.5
$ws.Range("D11").Value = 110.78
$ws.Range("E11").Value = 916.5
$ws.Range("G11").Value = 130.9
$ws.Range("N11").Value = 0.184
$ws.Range("O11").Value = 0.17899999999999999

$ws.Range("B12").Value = 1128
$ws.Range("D12").Value = 161.114
$ws.Range("E12").Value = 1057.5
$ws.Range("G12").Value = 151.07
$ws.Range("N12").Value = 0.187
$ws.Range("O12").Value = 0.17699999999999999

$ws.Range("B13").Value = 1057.5
$ws.Range("D13").Value = 151.07
$ws.Range("E13").Value = 1057.5
$ws.Range("G13").Value = 151.07
$ws.Range("N13").Value = 0.19600000000000001
$ws.Range("O13").Value = 0.16800000000000001

$ws.Range("B14").Value = 1128
$ws.Range("D14").Value = 161.114
# E14 also picks up the "thin side borders only" format (same as N14's).
$ws.Range("N14").Copy()
$ws.Range("E14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E14").Value = 846
$ws.Range("G14").Value = 120.8
$ws.Range("N14").Value = 0.16200000000000001
$ws.Range("O14").Value = 0.20200000000000001

$ws.Range("B15").Value = 916.5
$ws.Range("D15").Value = 130.9
$ws.Range("E15").Value = 916.5
$ws.Range("G15").Value = 130.9
$ws.Range("N15").Value = 0.17499999999999999
$ws.Range("O15").Value = 0.188

$ws.Range("B16").Value = 916.5
$ws.Range("D16").Value = 130.9
$ws.Range("E16").Value = 1057.5
$ws.Range("G16").Value = 151.07
$ws.Range("N16").Value = 0.21099999999999999
$ws.Range("O16").Value = 0.154

# ---- View state: selection + zoom, matching the author's last save ----
$ws.Range("L1").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 64
